$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.662.44"
$ws.Range("E2").Value = "  -3.16%  "

$ws.Range("D3").Value = "2.566.96"
$ws.Range("E3").Value = "  -5.69%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.56%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.85%  "

$ws.Range("E9").Value = "  -3.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.364"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.66%  "

$ws.Range("D13").Value = "3.025.79"
$ws.Range("E13").Value = "  -5.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.53%  "

$ws.Range("D15").Value = "61.576.57"
$ws.Range("E15").Value = "  -3.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000143"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.70%  "

$ws.Range("D17").Value = "2.576.58"
$ws.Range("E17").Value = "  -5.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.493"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.71%  "

$ws.Range("E25").Value = "  -1.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.95%  "

$ws.Range("D29").Value = "0.0₃0835"
$ws.Range("E29").Value = "  -5.64%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.86%  "

$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "159.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.50%  "

$ws.Range("E37").Value = "  +0.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "333.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.935"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.40%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.21%  "

$ws.Range("D45").Value = "2.136.25"
$ws.Range("E45").Value = "  +1.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.604"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.52%  "

$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.19%  "

$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0546"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0964"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0239"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.13%  "

